$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.440.28"
$ws.Range("E2").Value = "  +1.94%  "

$ws.Range("D3").Value = "'3.904.99"
$ws.Range("E3").Value = "  +0.40%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").Value = "'527.94"
$ws.Range("E5").Value = "  +9.08%  "

$ws.Range("D6").Value = "'144.58"
$ws.Range("E6").Value = "  -0.59%  "

$ws.Range("D7").Value = "'0.612"
$ws.Range("E7").Value = "  -1.55%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").Value = "'0.721"
$ws.Range("E9").Value = "  -2.53%  "

$ws.Range("E10").Value = "  -3.53%  "

$ws.Range("E11").Value = "  -6.30%  "

$ws.Range("D12").Value = "'42.22"
$ws.Range("E12").Value = "  -1.63%  "

$ws.Range("D13").Value = "'4.537.11"
$ws.Range("E13").Value = "  +0.64%  "

$ws.Range("D14").Value = "'10.26"
$ws.Range("E14").Value = "  -3.08%  "

$ws.Range("D15").Value = "'3.916.20"
$ws.Range("E15").Value = "  +0.37%  "

$ws.Range("D16").Value = "'14.05"
$ws.Range("E16").Value = "  -1.80%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.135"
$ws.Range("E17").Value = "  -0.72%  "

$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").Value = "'1.22"
$ws.Range("E18").Value = "  +7.27%  "

$ws.Range("D19").Value = "'19.75"

$ws.Range("D20").Value = "'69.410.73"
$ws.Range("E20").Value = "  +1.91%  "

$ws.Range("D21").Value = "'430.44"
$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("D22").Value = "'3.38"
$ws.Range("E22").Value = "  -5.26%  "

$ws.Range("D23").Value = "'14.26"
$ws.Range("E23").Value = "  -4.14%  "

$ws.Range("D24").Value = "'88.59"
$ws.Range("E24").Value = "  -0.98%  "

$ws.Range("E25").Value = "  +8.11%  "

$ws.Range("D26").Value = "'11.52"
$ws.Range("E26").Value = "  -2.36%  "

$ws.Range("D27").Value = "'10.62"
$ws.Range("E27").Value = "  -3.96%  "

$ws.Range("D28").Value = "'36.39"
$ws.Range("E28").Value = "  -2.95%  "

$ws.Range("D29").Value = "'692.07"
$ws.Range("E29").Value = "  -3.56%  "

$ws.Range("D30").Value = "'13.16"
$ws.Range("E30").Value = "  -3.44%  "

$ws.Range("E31").Value = "  -2.87%  "

$ws.Range("E32").Value = "  -3.51%  "

$ws.Range("D33").Value = "'67.68"
$ws.Range("E33").Value = "  +11.48%  "

$ws.Range("D34").Value = "'0.441"
$ws.Range("E34").Value = "  +10.73%  "

$ws.Range("D35").Value = "'5.96"
$ws.Range("E35").Value = "  -1.87%  "

$ws.Range("D36").Value = "'0.0₃0848"
$ws.Range("E36").Value = "  -4.28%  "

$ws.Range("D37").Value = "'39.89"
$ws.Range("E37").Value = "  -3.55%  "

$ws.Range("E38").Value = "  +3.46%  "

$ws.Range("D39").Value = "'0.997"
$ws.Range("E39").Value = "  +0.20%  "

$ws.Range("E40").Value = "  -0.12%  "

$ws.Range("D41").Value = "'0.0481"
$ws.Range("E41").Value = "  -3.09%  "

$ws.Range("D42").Value = "'3.13"
$ws.Range("E42").Value = "  +1.05%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'3.10"
$ws.Range("E43").Value = "  +3.96%  "

$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "'2.81"
$ws.Range("E44").Value = "  -6.39%  "

$ws.Range("D45").Value = "'3.36"
$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("D46").Value = "'0.140"
$ws.Range("E46").Value = "  -2.20%  "

$ws.Range("D47").Value = "'3.01"
$ws.Range("E47").Value = "  +7.19%  "

$ws.Range("D48").Value = "'0.0₆0351"
$ws.Range("E48").Value = "  +7.67%  "

$ws.Range("D49").Value = "'2.731.59"
$ws.Range("E49").Value = "  +10.83%  "

$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'144.90"
$ws.Range("E50").Value = "  -0.20%  "

$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").Value = "'3.28"
$ws.Range("E51").Value = "  -3.29%  "
